# Apply the Feb 28 2024 "cryptos" price/volume refresh from the upstream
# GitHub Actions scraper. Each row is the live coin table (B=Coin,
# C=Link, D=Price, E=Volume(1h)); two pairs of rows (Cosmos/Injective
# and NEARProtocol/ARBITRUM) swapped rank order, so those rows get new
# B/C/D/E content entirely instead of just refreshed numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '56.808.94'
$ws.Range("E2").Value = '  +0.77%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.239.99'
$ws.Range("E3").Value = '  +0.30%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.07%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '396.41'
$ws.Range("E5").Value = '  -0.53%  '

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.91'
$ws.Range("E6").Value = '  -2.72%  '

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.583'
$ws.Range("E7").Value = '  +5.41%  '

# Row 8: USDC
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.05%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -1.08%  '

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.25'
$ws.Range("E10").Value = '  -1.37%  '

# Row 12: TRON
$ws.Range("E12").Value = '  +2.00%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.747.31'
$ws.Range("E13").Value = '  +0.14%  '

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.32'
$ws.Range("E14").Value = '  +3.02%  '

# Row 15: Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.91'
$ws.Range("E15").Value = '  -1.17%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '3.238.78'
$ws.Range("E16").Value = '  +0.12%  '

# Row 17: Polygon
$ws.Range("E17").Value = '  -3.74%  '

# Row 18: Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.19'
$ws.Range("E18").Value = '  +6.76%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '56.698.99'
$ws.Range("E19").Value = '  +0.99%  '

# Row 20: ImmutableX
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  -1.11%  '

# Row 21: ShibaInu
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000111'
$ws.Range("E21").Value = '  +10.25%  '

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.02'
$ws.Range("E22").Value = '  -0.77%  '

# Row 23: BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '291.08'
$ws.Range("E23").Value = '  +0.71%  '

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.50'
$ws.Range("E24").Value = '  +0.61%  '

# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.17'
$ws.Range("E25").Value = '  -1.82%  '

# Row 26: EthereumClassic
$ws.Range("E26").Value = '  -0.39%  '

# Row 27: LEO
$ws.Range("E27").Value = '  -0.78%  '

# Row 28: Filecoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.78'
$ws.Range("E28").Value = '  -4.92%  '

# Row 29: Kaspa
$ws.Range("E29").Value = '  -1.57%  '

# Row 30: RenderToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("E30").Value = '  -3.86%  '

# Row 31: Dai
$ws.Range("E31").Value = '  +0.05%  '

# Row 32: InjectiveProtocol
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '41.62'
$ws.Range("E32").Value = '  +12.38%  '

# Row 33: Cosmos
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.18'
$ws.Range("E33").Value = '  -0.87%  '

# Row 34: Hedera
$ws.Range("E34").Value = '  -3.04%  '

# Row 35: VeChain
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0484'
$ws.Range("E35").Value = '  -3.48%  '

# Row 36: Toncoin
$ws.Range("E36").Value = '  +1.44%  '

# Row 37: OKB
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.33'
$ws.Range("E37").Value = '  +0.46%  '

# Row 38: FirstDigitalUSD
$ws.Range("E38").Value = '  -0.20%  '

# Row 39: LidoDAOToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.47'
$ws.Range("E39").Value = '  -3.87%  '

# Row 40: Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -4.08%  '

# Row 41: Monero
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '136.65'
$ws.Range("E41").Value = '  -1.68%  '

# Row 42: Stellar
$ws.Range("E42").Value = '  +3.51%  '

# Row 43: ARBITRUM
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.88'
$ws.Range("E43").Value = '  -2.67%  '

# Row 44: NEARProtocol
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.94'
$ws.Range("E44").Value = '  -2.75%  '

# Row 45: Celestia
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.85'
$ws.Range("E45").Value = '  -0.77%  '

# Row 46: TheGraph
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.276'
$ws.Range("E46").Value = '  -2.82%  '

# Row 47: WEMIXToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +8.48%  '

# Row 48: EnergySwap
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.54'
$ws.Range("E48").Value = '  -0.01%  '

# Row 49: Maker
$ws.Range("D49").Value = '2.149.72'
$ws.Range("E49").Value = '  +1.15%  '

# Row 50: ApeXProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("E50").Value = '  -5.79%  '

# Row 51: ThetaToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.95'
$ws.Range("E51").Value = '  -5.73%  '
